{"js": "// Replace the date and each of the division problems with their new values.\n// Every old string is unique within the document, so a targeted\n// search-and-replace (which preserves run formatting) is safe here.\nconst replacements = [\n  [\"2025-05-02 Friday\", \"2025-05-03 Saturday\"],\n  [\"561\u00f74=\", \"319\u00f76=\"],\n  [\"872\u00f74=\", \"233\u00f72=\"],\n  [\"850\u00f73=\", \"883\u00f73=\"],\n  [\"469\u00f78=\", \"571\u00f78=\"],\n  [\"505\u00f73=\", \"163\u00f77=\"],\n  [\"123\u00f75=\", \"223\u00f77=\"],\n  [\"764\u00f77=\", \"300\u00f73=\"],\n  [\"220\u00f77=\", \"511\u00f77=\"],\n  [\"925\u00f76=\", \"220\u00f78=\"],\n  [\"775\u00f79=\", \"347\u00f77=\"],\n  [\"499\u00f74=\", \"855\u00f77=\"],\n  [\"570\u00f72=\", \"471\u00f72=\"],\n  [\"143\u00f73=\", \"195\u00f74=\"],\n  [\"264\u00f79=\", \"804\u00f76=\"],\n  [\"233\u00f73=\", \"105\u00f72=\"],\n  [\"726\u00f75=\", \"670\u00f77=\"],\n  [\"210\u00f74=\", \"737\u00f73=\"],\n  [\"436\u00f75=\", \"151\u00f74=\"],\n  [\"479\u00f78=\", \"750\u00f75=\"],\n  [\"736\u00f72=\", \"417\u00f76=\"],\n  [\"341\u00f78=\", \"349\u00f72=\"],\n  [\"259\u00f75=\", \"567\u00f78=\"],\n  [\"137\u00f79=\", \"400\u00f79=\"],\n  [\"574\u00f74=\", \"111\u00f78=\"],\n  [\"910\u00f76=\", \"233\u00f74=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the date and each of the division problems with their new values.\n# Every old string is unique within the document, so Find/Replace (which\n# preserves run formatting) is safe here.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"2025-05-02 Friday\", \"2025-05-03 Saturday\"),\n    @(\"561\u00f74=\", \"319\u00f76=\"),\n    @(\"872\u00f74=\", \"233\u00f72=\"),\n    @(\"850\u00f73=\", \"883\u00f73=\"),\n    @(\"469\u00f78=\", \"571\u00f78=\"),\n    @(\"505\u00f73=\", \"163\u00f77=\"),\n    @(\"123\u00f75=\", \"223\u00f77=\"),\n    @(\"764\u00f77=\", \"300\u00f73=\"),\n    @(\"220\u00f77=\", \"511\u00f77=\"),\n    @(\"925\u00f76=\", \"220\u00f78=\"),\n    @(\"775\u00f79=\", \"347\u00f77=\"),\n    @(\"499\u00f74=\", \"855\u00f77=\"),\n    @(\"570\u00f72=\", \"471\u00f72=\"),\n    @(\"143\u00f73=\", \"195\u00f74=\"),\n    @(\"264\u00f79=\", \"804\u00f76=\"),\n    @(\"233\u00f73=\", \"105\u00f72=\"),\n    @(\"726\u00f75=\", \"670\u00f77=\"),\n    @(\"210\u00f74=\", \"737\u00f73=\"),\n    @(\"436\u00f75=\", \"151\u00f74=\"),\n    @(\"479\u00f78=\", \"750\u00f75=\"),\n    @(\"736\u00f72=\", \"417\u00f76=\"),\n    @(\"341\u00f78=\", \"349\u00f72=\"),\n    @(\"259\u00f75=\", \"567\u00f78=\"),\n    @(\"137\u00f79=\", \"400\u00f79=\"),\n    @(\"574\u00f74=\", \"111\u00f78=\"),\n    @(\"910\u00f76=\", \"233\u00f74=\")\n)\n\nforeach ($pair in $pairs) {\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($pair[0], $false, $false, $false, $false, $false, $true, 1, $false, $pair[1], 2)\n}\n"}
